$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A for the "Lab. #" identifier - shifts the existing
# ratio columns (old A:T) right to (new B:U).
$ws.Columns("A:A").Insert()
$ws.Range("A1").Value = "Lab. #"
$ws.Columns("A:A").ColumnWidth = 6.75

# Fill in the lab numbers for each data row (rows 2-16).
$labNums = @(9186, 9715, 9186, 9716, 9186, 9717, 9186, 9718, 9186, 9719, 9186, 9720, 9186, 9721, 9186)
for ($i = 0; $i -lt $labNums.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labNums[$i]
}

# Highlight the repeated/duplicate-standard rows (Lab. # 9186) with a light
# green fill so they stand out from the unique sample rows.
$highlightColor = 12379352
$highlightRows = @(2, 4, 6, 8, 10, 12, 14, 16)
foreach ($r in $highlightRows) {
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 21)).Interior.Color = $highlightColor
}
